$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$halfPi = 1.5707963267948966

for ($r = 1; $r -le 51; $r++) {
    $bCell = $ws.Cells.Item($r, 2)  # Column B
    $cCell = $ws.Cells.Item($r, 3)  # Column C
    $dCell = $ws.Cells.Item($r, 4)  # Column D
    $eCell = $ws.Cells.Item($r, 5)  # Column E
    $fCell = $ws.Cells.Item($r, 6)  # Column F

    $bOld = $bCell.Value()
    $cOld = $cCell.Value()
    $dOld = $dCell.Value()
    $eOld = $eCell.Value()
    $fOld = $fCell.Value()

    $bCell.Value = -1 * $bOld
    $cCell.Value = $halfPi - $cOld
    $dCell.Value = -1 * $dOld
    $eCell.Value = $halfPi - $eOld
    $fCell.Value = -1 * $fOld
}
